$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 4.1
$ws.Range("I2").Value = 3.9
$ws.Range("J2").Value = 2.38
$ws.Range("Q2").Value = 1.62
$ws.Range("R2").Value = 2.25
$ws.Range("U2").Value = 1.62
$ws.Range("V2").Value = 2.2
$ws.Range("W2").Value = 9
$ws.Range("AD2").Value = 8
$ws.Range("AO2").Value = 9
$ws.Range("AZ2").Value = 26
$ws.Range("BB2").Value = 81

# Row 9
$ws.Range("H9").Value = 3.2
$ws.Range("I9").Value = 3.75
$ws.Range("J9").Value = 2.88
$ws.Range("K9").Value = 1.95
$ws.Range("Q9").Value = 2.5
$ws.Range("R9").Value = 1.5
$ws.Range("X9").Value = 8.5
$ws.Range("AC9").Value = 7
$ws.Range("AJ9").Value = 13
$ws.Range("AL9").Value = 34
$ws.Range("AM9").Value = 41
$ws.Range("AO9").Value = 12
$ws.Range("AP9").Value = 26
$ws.Range("AV9").Value = 67

# Row 10
$ws.Range("Q10").Value = 1.73
$ws.Range("R10").Value = 2.08

# Row 11
$ws.Range("G11").Value = 5
$ws.Range("H11").Value = 3.45
$ws.Range("I11").Value = 1.65
$ws.Range("J11").Value = 5.4
$ws.Range("K11").Value = 2.1
$ws.Range("L11").Value = 2.22
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 6.5
$ws.Range("O11").Value = 1.38
$ws.Range("P11").Value = 2.82
$ws.Range("Q11").Value = 2.1
$ws.Range("R11").Value = 1.65
$ws.Range("S11").Value = 1.44
$ws.Range("T11").Value = 2.57
$ws.Range("U11").Value = 2.05
$ws.Range("V11").Value = 1.7
$ws.Range("W11").Value = 11.75
$ws.Range("X11").Value = 29
$ws.Range("Y11").Value = 17
$ws.Range("Z11").Value = 100
$ws.Range("AA11").Value = 60
$ws.Range("AB11").Value = 70
$ws.Range("AC11").Value = 6.5
$ws.Range("AD11").Value = 6.9
$ws.Range("AE11").Value = 19
$ws.Range("AF11").Value = 110
$ws.Range("AH11").Value = 5.6
$ws.Range("AI11").Value = 6.9
$ws.Range("AK11").Value = 12
$ws.Range("AL11").Value = 15
$ws.Range("AN11").Value = 6.7
$ws.Range("AO11").Value = 32
$ws.Range("AP11").Value = 40
$ws.Range("AQ11").Value = 200
$ws.Range("AR11").Value = 250
$ws.Range("AT11").Value = 2.57
$ws.Range("AU11").Value = 8.25
$ws.Range("AV11").Value = 90
$ws.Range("AX11").Value = 3.35
$ws.Range("AY11").Value = 8.25
$ws.Range("AZ11").Value = 20
$ws.Range("BA11").Value = 28
$ws.Range("BB11").Value = 70
